$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Name = "sdjcdi"

$ws.Cells.Item(1, 10).Value2 = 60.76769471168518
$ws.Cells.Item(2, 10).Value2 = 76.96686887741089
$ws.Cells.Item(3, 10).Value2 = 96.61346435546875
$ws.Cells.Item(4, 2).Value2 = 2592
$ws.Cells.Item(4, 5).Value2 = 27
$ws.Cells.Item(4, 8).Value2 = 98.95793130065611
$ws.Cells.Item(4, 9).Value2 = 0.01321928460342146
$ws.Cells.Item(4, 10).Value2 = 77.80708527565002
$ws.Cells.Item(5, 2).Value2 = 2022
$ws.Cells.Item(5, 4).Value2 = 2005
$ws.Cells.Item(5, 6).Value2 = 21
$ws.Cells.Item(5, 7).Value2 = 98.9634748272458
$ws.Cells.Item(5, 8).Value2 = 99.208312716477
$ws.Cells.Item(5, 9).Value2 = 0.01825357671435619
$ws.Cells.Item(5, 10).Value2 = 93.14046573638916
$ws.Cells.Item(6, 10).Value2 = 108.6741955280304
$ws.Cells.Item(7, 10).Value2 = 108.4913229942322
$ws.Cells.Item(8, 10).Value2 = 152.9010188579559
$ws.Cells.Item(9, 10).Value2 = 138.4324190616608
$ws.Cells.Item(10, 2).Value2 = 1893
$ws.Cells.Item(10, 4).Value2 = 1793
$ws.Cells.Item(10, 5).Value2 = 99
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 99.94425863991081
$ws.Cells.Item(10, 8).Value2 = 94.76744186046511
$ws.Cells.Item(10, 9).Value2 = 0.05571030640668524
$ws.Cells.Item(10, 10).Value2 = 163.4708733558655
$ws.Cells.Item(11, 10).Value2 = 181.3979864120483
$ws.Cells.Item(12, 10).Value2 = 205.4927563667297
$ws.Cells.Item(13, 10).Value2 = 248.336820602417
$ws.Cells.Item(14, 10).Value2 = 195.9933547973633
$ws.Cells.Item(15, 2).Value2 = 2281
$ws.Cells.Item(15, 5).Value2 = 3
$ws.Cells.Item(15, 8).Value2 = 99.86842105263158
$ws.Cells.Item(15, 9).Value2 = 0.00131694468832309
$ws.Cells.Item(15, 10).Value2 = 250.1790881156921
$ws.Cells.Item(16, 10).Value2 = 218.1064755916595
$ws.Cells.Item(17, 10).Value2 = 267.8450391292572
$ws.Cells.Item(18, 10).Value2 = 277.5529115200043
$ws.Cells.Item(19, 2).Value2 = 1519
$ws.Cells.Item(19, 5).Value2 = 1
$ws.Cells.Item(19, 8).Value2 = 99.93412384716733
$ws.Cells.Item(19, 9).Value2 = 0.0006587615283267457
$ws.Cells.Item(19, 10).Value2 = 251.3054463863373
$ws.Cells.Item(20, 10).Value2 = 331.1922447681427
$ws.Cells.Item(21, 10).Value2 = 239.7215094566345
$ws.Cells.Item(22, 10).Value2 = 269.7521753311157
$ws.Cells.Item(23, 2).Value2 = 2133
$ws.Cells.Item(23, 4).Value2 = 2130
$ws.Cells.Item(23, 6).Value2 = 5
$ws.Cells.Item(23, 7).Value2 = 99.76580796252928
$ws.Cells.Item(23, 8).Value2 = 99.90619136960601
$ws.Cells.Item(23, 9).Value2 = 0.003277153558052435
$ws.Cells.Item(23, 10).Value2 = 290.8733282089233
$ws.Cells.Item(24, 2).Value2 = 2923
$ws.Cells.Item(24, 4).Value2 = 2901
$ws.Cells.Item(24, 5).Value2 = 21
$ws.Cells.Item(24, 6).Value2 = 78
$ws.Cells.Item(24, 7).Value2 = 97.38167170191339
$ws.Cells.Item(24, 8).Value2 = 99.28131416837782
$ws.Cells.Item(24, 9).Value2 = 0.03322147651006711
$ws.Cells.Item(24, 10).Value2 = 280.3202874660492
$ws.Cells.Item(25, 10).Value2 = 194.4791741371155
$ws.Cells.Item(26, 2).Value2 = 1850
$ws.Cells.Item(26, 4).Value2 = 1846
$ws.Cells.Item(26, 6).Value2 = 13
$ws.Cells.Item(26, 7).Value2 = 99.30069930069931
$ws.Cells.Item(26, 8).Value2 = 99.83775013520822
$ws.Cells.Item(26, 9).Value2 = 0.008602150537634409
$ws.Cells.Item(26, 10).Value2 = 364.4736130237579
$ws.Cells.Item(27, 10).Value2 = 178.7861235141754
$ws.Cells.Item(28, 10).Value2 = 70.77923989295959
$ws.Cells.Item(29, 10).Value2 = 270.9186751842499
$ws.Cells.Item(30, 10).Value2 = 107.5567719936371
$ws.Cells.Item(31, 10).Value2 = 39.68196082115173
$ws.Cells.Item(32, 10).Value2 = 466.8973119258881
$ws.Cells.Item(33, 2).Value2 = 3361
$ws.Cells.Item(33, 4).Value2 = 3360
$ws.Cells.Item(33, 6).Value2 = 2
$ws.Cells.Item(33, 7).Value2 = 99.94051160023795
$ws.Cells.Item(33, 9).Value2 = 0.0005947071067499256
$ws.Cells.Item(33, 10).Value2 = 110.2967150211334
$ws.Cells.Item(34, 10).Value2 = 409.9186980724335
$ws.Cells.Item(35, 2).Value2 = 2048
$ws.Cells.Item(35, 4).Value2 = 2047
$ws.Cells.Item(35, 6).Value2 = 0
$ws.Cells.Item(35, 7).Value2 = 100
$ws.Cells.Item(35, 9).Value2 = 0
$ws.Cells.Item(35, 10).Value2 = 513.3081405162811
$ws.Cells.Item(36, 10).Value2 = 325.8544006347656
$ws.Cells.Item(37, 2).Value2 = 2347
$ws.Cells.Item(37, 4).Value2 = 2346
$ws.Cells.Item(37, 6).Value2 = 136
$ws.Cells.Item(37, 7).Value2 = 94.52054794520548
$ws.Cells.Item(37, 9).Value2 = 0.05477245267821184
$ws.Cells.Item(37, 10).Value2 = 431.7027928829193
$ws.Cells.Item(38, 10).Value2 = 314.4228749275208
$ws.Cells.Item(39, 2).Value2 = 2055
$ws.Cells.Item(39, 5).Value2 = 6
$ws.Cells.Item(39, 8).Value2 = 99.70788704965921
$ws.Cells.Item(39, 9).Value2 = 0.004870920603994155
$ws.Cells.Item(39, 10).Value2 = 455.6879987716675
$ws.Cells.Item(40, 10).Value2 = 529.5843977928162
$ws.Cells.Item(41, 10).Value2 = 442.6577336788177
$ws.Cells.Item(42, 2).Value2 = 1780
$ws.Cells.Item(42, 4).Value2 = 1778
$ws.Cells.Item(42, 6).Value2 = 1
$ws.Cells.Item(42, 7).Value2 = 99.94378864530636
$ws.Cells.Item(42, 8).Value2 = 99.94378864530636
$ws.Cells.Item(42, 9).Value2 = 0.001123595505617978
$ws.Cells.Item(42, 10).Value2 = 444.8405418395996
$ws.Cells.Item(43, 10).Value2 = 352.9401865005493
$ws.Cells.Item(44, 10).Value2 = 34.61887836456299
